# Daily attendance processing - 2025-10-31 10:48:40
# Reorders the "Recorded By" (column G) list on the active sheet so that
# any "System"/"system" entries are moved to the front of the
# comma-separated list, preserving the relative order of the remaining
# entries. Rows that already start with "System"/"system", rows with no
# "System" entry at all, and empty cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    $parts = $val -split ", "

    if ($parts.Count -le 1) {
        continue
    }

    $systemParts = @()
    $otherParts = @()

    foreach ($p in $parts) {
        if ($p -eq "System" -or $p -eq "system") {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Count -eq 0) {
        continue
    }

    if ($parts[0] -eq "System" -or $parts[0] -eq "system") {
        continue
    }

    $newParts = $systemParts + $otherParts
    $newVal = $newParts -join ", "

    $cell.Value = $newVal
}
